# Automatische test-sync: 2025-06-27 22:19:50
# Append a new "Logs" row for the test mail "Wanneer zijn jullie open?"
# and refresh the "Dashboard" summary sheet to reflect the new counts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 5 with the new test-mail entry.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A5").Value = "Wanneer zijn jullie open?"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Testmail #1: Wanneer zijn jullie open?"
$logs.Range("D5").Value = "Openingstijden / Locatie"
$logs.Range("E5").Value = "Beste klant,`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F5").Value = "2025-06-27 22:19:31"
$logs.Range("G5").Value = "Ja"
$logs.Range("H5").Value = "Nee"
$logs.Range("I5").Value = "Ja"

# Extend the conditional-formatting ranges so the new row is covered too
# (D/G/H/I columns), mirroring the widened dimension A1:I4 -> A1:I5.
$logs.Range("D2:D4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D5"))
$logs.Range("G2:G4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G5"))
$logs.Range("H2:H4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H5"))
$logs.Range("I2:I4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I5"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: re-sync the category counts now that
#    "Openingstijden / Locatie" occurs twice, reordering the rows.
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A2").Value = "Openingstijden / Locatie"
$dashboard.Range("B2").Value = 2
$dashboard.Range("A3").Value = "Planning / Afspraak"
$dashboard.Range("B3").Value = 1
$dashboard.Range("A4").Value = "Overig"
$dashboard.Range("B4").Value = 1
